$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '42.880.73'
$ws.Range("E2").Value = '  +0.24%  '

$ws.Range("D3").Value = '2.534.97'
$ws.Range("E3").Value = '  +0.46%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  +0.01%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '318.29'
$ws.Range("E5").Value = '  +4.37%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '95.83'
$ws.Range("E6").Value = '  -1.41%  '

$ws.Range("E7").Value = '  +0.08%  '

$ws.Range("E8").Value = '  -0.07%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.534'
$ws.Range("E9").Value = '  -1.49%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '36.34'
$ws.Range("E10").Value = '  -1.33%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0812'
$ws.Range("E11").Value = '  -0.24%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '7.64'
$ws.Range("E12").Value = '  +1.31%  '

$ws.Range("E13").Value = '  -0.46%  '

$ws.Range("D14").Value = '2.923.35'
$ws.Range("E14").Value = '  +0.58%  '

$ws.Range("D15").Value = '2.549.54'
$ws.Range("E15").Value = '  +2.12%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '15.34'
$ws.Range("E16").Value = '  +1.64%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.849'
$ws.Range("E17").Value = '  -1.57%  '

$ws.Range("D18").Value = '42.937.93'
$ws.Range("E18").Value = '  +0.39%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '13.07'
$ws.Range("E19").Value = '  +0.62%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.65'
$ws.Range("E20").Value = '  +2.82%  '

$ws.Range("E21").Value = '  -0.66%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '70.19'
$ws.Range("E22").Value = '  -1.47%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '252.11'
$ws.Range("E23").Value = '  +0.17%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.97'
$ws.Range("E24").Value = '  +1.83%  '

$ws.Range("E25").Value = '  -0.32%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '27.07'
$ws.Range("E26").Value = '  +0.37%  '

$ws.Range("E27").Value = '  +0.08%  '

$ws.Range("E28").Value = '  +3.79%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '39.89'
$ws.Range("E29").Value = '  +4.68%  '

$ws.Range("E30").Value = '  -0.89%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.08'
$ws.Range("E31").Value = '  +1.59%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '153.83'
$ws.Range("E32").Value = '  -1.78%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '2.13'

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.33'
$ws.Range("E34").Value = '  +1.64%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '19.03'
$ws.Range("E35").Value = '  +3.35%  '

$ws.Range("E36").Value = '  -0.28%  '

$ws.Range("E37").Value = '  +0.06%  '

$ws.Range("E38").Value = '  -3.70%  '

$ws.Range("E39").Value = '  -0.36%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '23.68'
$ws.Range("E40").Value = '  -1.16%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.29'
$ws.Range("E41").Value = '  +9.59%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.81'
$ws.Range("E42").Value = '  -1.02%  '

$ws.Range("E43").Value = '  +1.18%  '

$ws.Range("B44").Value = 'NEARProtocol'
$ws.Range("C44").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '3.33'
$ws.Range("E44").Value = '  -2.07%  '

$ws.Range("B45").Value = 'FirstDigitalUSD'
$ws.Range("C45").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.00'
$ws.Range("E45").Value = '  +0.40%  '

$ws.Range("D46").Value = '2.017.87'
$ws.Range("E46").Value = '  -0.69%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '85.75'
$ws.Range("E47").Value = '  +0.27%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '8.78'
$ws.Range("E48").Value = '  -2.07%  '

$ws.Range("D49").Value = '2.781.20'
$ws.Range("E49").Value = '  +0.52%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '74.30'
$ws.Range("E50").Value = '  +2.67%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '102.47'
$ws.Range("E51").Value = '  +0.34%  '
